$d = $word.ActiveDocument
$t = $d.Tables(1)
$newValues = @("81-3=","36+57=","65-59=","43-19=","39-8=","43+35=","84-21=","43+17=","74+0=","61+35=","55+8=","84-30=","38-37=","85-0=","90-85=","61-54=","81-30=","27+66=","81-43=","71+23=","69-2=","18+67=","62-39=","5+3=","63+21=","27+57=","9+45=","8+29=","19+30=","32-5=","29-20=","35+37=","68+9=","92-38=","1+82=","75-26=","0+7=","52-18=","42+32=","18+54=","19-2=","46+10=","38+57=","60+11=","25+1=","88+0=","36-10=","29+58=","32+53=","53-29=","98-33=","57-56=","61+28=","37-2=","12+10=","55+39=","53-48=","44+2=","24+43=","67-63=","71-63=","44-30=","17+4=","72-42=","77-41=","88-67=","83-60=","4+23=","77-68=","35-8=","43+1=","30+11=","24+68=","28+11=","0+77=","33+28=","62-32=","15+28=","60-11=","97-50=","36+26=","8+11=","27+22=","22-5=","30+1=","13+50=","5+41=","88-85=","78-27=","83-69=","41-13=","68+16=","49+7=","62-19=","13-9=","76-31=","25+15=","74-47=","77+21=","98-56=")
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}
Write-Host "Updated" $idx "cells"
